# Commit: "add default hit time"
# Adds a new schema row (31) describing a "DefaultHitTime" field to the
# Property sheet of Skill.xlsx, following the same column layout as the
# other rows (A=Id, B=Type, C=Public, D=Private, E=Save, F=View,
# G=Index, H=SaveInterval, I=RelationValue, J=Desc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "DefaultHitTime"
$ws.Range("B31").Value = "float"
$ws.Range("C31").Value = $false
$ws.Range("D31").Value = $false
$ws.Range("E31").Value = $false
$ws.Range("F31").Value = $true
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = "Friend"
$ws.Range("J31").Value = "缺省打击时间（本来应该打到但是物理没碰撞到或者其他原因）"

# The Id/Type/RelationValue/Desc columns are formatted as Text in every
# other schema row (A,B,I,J) - match that here. C:H stay General, same as
# every other row.
$ws.Range("A31").NumberFormat = "@"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("I31").NumberFormat = "@"
$ws.Range("J31").NumberFormat = "@"

# Leave the selection just past the newly-added row, mirroring where the
# cursor was left after the edit.
[void]$ws.Range("J32").Select()
